$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.918.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.26%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.551.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.20%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.55%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'206.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.27%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.488"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.42%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.53%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'22.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.85%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D10").Value = "'0.0595"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.96%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.70%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.772.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.24%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.545.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.60%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.86%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.81%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'26.904.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.25%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'61.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.58%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +2.62%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'217.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.65%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.45%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.51%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.95%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.18%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.27%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'153.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.75%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.39%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.46%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.66%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.41%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.29%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.14%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.31%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'Maker"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'1.418.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.93%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'InternetComputer(DFINITY)"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'3.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.47%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'TrustWalletToken"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'1.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +13.87%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'LidoDAOToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'1.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.78%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.28%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.26%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.87%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.39%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.54%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.35%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.77%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.60%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'64.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.90%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.56%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.686.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.26%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'87.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.00%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +1.47%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0₆0101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.67%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.09%  "
$ws.Range("E51").Style = "Normal"
